$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("2010年", 16395.7325, 613.1541999999999, 17814.0904, 920.3656999999999, 9244.3326, 72.74290000000001, 3181.7496, 17.2351),
    @("2011年", 19155.5027, 523.4813, 25706.5729, 868.3763, 12327.0766, 54.7509, 3150.8089, 10.3162),
    @("2012年", 19743.9093, 418.0822, 33210.7617, 763.0562, 14711.402, 45.2615, 2903.7823, 7.5317),
    @("2013年", 22144.3942, 378.1633, 41574.138, 698.0399, 18419.6454, 32.5938, 2761.3338, 5.0742),
    @("2014年", 21690.7331, 266.1907, 49820.2606, 550.5103, 20404.8022, 18.5819, 2282.0775, 2.4581),
    @("2015年", 19381.5699, 195.0928, 53610.8655, 412.6776, 20561.1417, 11.4432, 1804.8695, 1.1856),
    @("2016年", 20369.3927, 148.5083, 58584.4087, 354.7742, 21381.3694, 8.5746, 1732.5383, 1.0466),
    @("2017年", 20088.6071, 132.707, 65201.1241, 318.9941, 22460.2189, 6.8685, 1589.5544, 0.4547),
    @("2018年", 20088.534, 119.9816, 74832.361, 275.3736, 23436.6137, 6.9767, 1503.2988, 0.3752),
    @("2019年", 20622.6004, 116.8259, 84092.9602, 277.9759, 25659.2087, 6.0372, 1418.2573, 0.3993),
    @("2020年", 19925.4001, 100.2176, 93127.8135, 247.776, 26793.2174, 5.9577, 1242.2307, 0.3325),
    @("2021年", 19365.2765, 93.5351, 99715.901, 226.5763, 27073.3581, 4.7739, 1122.4223, 0.2335)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    for ($c = 1; $c -lt $vals.Length; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = [double]$vals[$c]
    }
}

# Remove the now-stale rows 14-18 (old data shifted out of range)
$ws.Range("A14:I18").EntireRow.Delete() | Out-Null
